# Auto-generated script applying scheduled price-data refresh to Sheets workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2486.6155
$ws.Range("J17").Value = 2486.6155
$ws.Range("L17").Value = 7459.8465
$ws.Range("N17").Value = -7795.8465
$ws.Range("H116").Value = 13217318
$ws.Range("I116").Value = 17479378
$ws.Range("J116").Value = 4928.8
$ws.Range("K116").Value = 17479378
$ws.Range("L116").Value = 4928.8
$ws.Range("M116").Value = -17475936
$ws.Range("N116").Value = -11812.8
$ws.Range("H129").Value = 1840.8572
$ws.Range("J129").Value = 2930.6667
$ws.Range("L129").Value = 8792.000100000001
$ws.Range("N129").Value = -18792.0001
$ws.Range("H131").Value = 3288.7273
$ws.Range("J131").Value = 5197.25
$ws.Range("L131").Value = 15591.75
$ws.Range("N131").Value = -25671.75
$ws.Range("H133").Value = 99774.5
$ws.Range("J133").Value = 99774.5
$ws.Range("L133").Value = 99774.5
$ws.Range("N133").Value = -109894.5
$ws.Range("H138").Value = 7498.9834
$ws.Range("I138").Value = 3552.3333
$ws.Range("K138").Value = 10656.9999
$ws.Range("M138").Value = -5516.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9906806
$ws.Range("I2").Value = 873897.2
$ws.Range("J2").Value = 35715116
$ws.Range("K2").Value = 873897.2
$ws.Range("L2").Value = 35715116
$ws.Range("M2").Value = -873784.2
$ws.Range("N2").Value = -35715342
$ws.Range("H31").Value = 2287.25
$ws.Range("I31").Value = 2287.25
$ws.Range("K31").Value = 2287.25
$ws.Range("M31").Value = -1993.25
$ws.Range("H32").Value = 13295.214
$ws.Range("I32").Value = 12693.302
$ws.Range("K32").Value = 12693.302
$ws.Range("M32").Value = -12406.302
$ws.Range("H45").Value = 3808.0588
$ws.Range("J45").Value = 4061.1428
$ws.Range("L45").Value = 4061.1428
$ws.Range("N45").Value = -4815.1428
$ws.Range("H61").Value = 4017.2778
$ws.Range("I61").Value = 3305
$ws.Range("K61").Value = 3305
$ws.Range("M61").Value = -3093
$ws.Range("H74").Value = 14707790
$ws.Range("I74").Value = 20834578
$ws.Range("J74").Value = 3499.4
$ws.Range("K74").Value = 20834578
$ws.Range("L74").Value = 3499.4
$ws.Range("M74").Value = -20833704
$ws.Range("N74").Value = -5247.4
$ws.Range("H77").Value = 14707790
$ws.Range("I77").Value = 20834578
$ws.Range("J77").Value = 3499.4
$ws.Range("K77").Value = 104172890
$ws.Range("L77").Value = 17497
$ws.Range("M77").Value = -104168522
$ws.Range("N77").Value = -26233
$ws.Range("H97").Value = 1072.7059
$ws.Range("I97").Value = 468.92307
$ws.Range("K97").Value = 468.92307
$ws.Range("M97").Value = 27.07693
$ws.Range("H116").Value = 9906806
$ws.Range("I116").Value = 873897.2
$ws.Range("J116").Value = 35715116
$ws.Range("K116").Value = 873897.2
$ws.Range("L116").Value = 35715116
$ws.Range("M116").Value = -871603.2
$ws.Range("N116").Value = -35719704
$ws.Range("H122").Value = 8268.25
$ws.Range("I122").Value = 6429.7
$ws.Range("K122").Value = 19289.1
$ws.Range("M122").Value = -16839.1
$ws.Range("H132").Value = 15015.745
$ws.Range("I132").Value = 16721.324
$ws.Range("K132").Value = 50163.972
$ws.Range("M132").Value = -47633.972
$ws.Range("H136").Value = 4017.2778
$ws.Range("I136").Value = 3305
$ws.Range("K136").Value = 9915
$ws.Range("M136").Value = -7365

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9906806
$ws.Range("I3").Value = 873897.2
$ws.Range("J3").Value = 35715116
$ws.Range("K3").Value = 873897.2
$ws.Range("L3").Value = 35715116
$ws.Range("M3").Value = -873783.2
$ws.Range("N3").Value = -35715344
$ws.Range("H99").Value = 92804510
$ws.Range("I99").Value = 92804510
$ws.Range("K99").Value = 92804510
$ws.Range("M99").Value = -92803012
$ws.Range("H105").Value = 44119400
$ws.Range("I105").Value = 46876776
$ws.Range("K105").Value = 46876776
$ws.Range("M105").Value = -46875029
$ws.Range("H134").Value = 2069.9688
$ws.Range("I134").Value = 1901.5385
$ws.Range("J134").Value = 2799.8333
$ws.Range("K134").Value = 5704.6155
$ws.Range("L134").Value = 8399.499899999999
$ws.Range("M134").Value = -3169.6155
$ws.Range("N134").Value = -13469.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1724
$ws.Range("I2").Value = 600.5
$ws.Range("J2").Value = 2847.5
$ws.Range("K2").Value = 600.5
$ws.Range("L2").Value = 2847.5
$ws.Range("M2").Value = -487.5
$ws.Range("N2").Value = -3073.5
$ws.Range("H16").Value = 1927.8889
$ws.Range("I16").Value = 1956.375
$ws.Range("J16").Value = 1700
$ws.Range("K16").Value = 1956.375
$ws.Range("L16").Value = 1700
$ws.Range("M16").Value = -1669.375
$ws.Range("N16").Value = -2274
$ws.Range("H31").Value = 31253328
$ws.Range("I31").Value = 34485396
$ws.Range("K31").Value = 34485396
$ws.Range("M31").Value = -34485101
$ws.Range("H34").Value = 31253328
$ws.Range("I34").Value = 34485396
$ws.Range("K34").Value = 34485396
$ws.Range("M34").Value = -34485194
$ws.Range("H74").Value = 79992.664
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 79992.664
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 79992.664
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -81740.664
$ws.Range("H77").Value = 79992.664
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 79992.664
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 239977.992
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -248713.992
$ws.Range("H86").Value = 10610.483
$ws.Range("I86").Value = 11266
$ws.Range("J86").Value = 10137.056
$ws.Range("K86").Value = 11266
$ws.Range("L86").Value = 10137.056
$ws.Range("M86").Value = -10143
$ws.Range("N86").Value = -12383.056
$ws.Range("H89").Value = 10610.483
$ws.Range("I89").Value = 11266
$ws.Range("J89").Value = 10137.056
$ws.Range("K89").Value = 56330
$ws.Range("L89").Value = 50685.28
$ws.Range("M89").Value = -50714
$ws.Range("N89").Value = -61917.28
$ws.Range("H113").Value = 1927.8889
$ws.Range("I113").Value = 1956.375
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 1956.375
$ws.Range("L113").Value = 1700
$ws.Range("M113").Value = 213.625
$ws.Range("N113").Value = -6040
$ws.Range("H132").Value = 18527746
$ws.Range("I132").Value = 21864648
$ws.Range("K132").Value = 65593944
$ws.Range("M132").Value = -65591414
$ws.Range("H134").Value = 1742.7941
$ws.Range("J134").Value = 1206.75
$ws.Range("L134").Value = 3620.25
$ws.Range("N134").Value = -8690.25
$ws.Range("H141").Value = 127386.55
$ws.Range("J141").Value = 127386.55
$ws.Range("L141").Value = 127386.55
$ws.Range("N141").Value = -137746.55

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3440.2778
$ws.Range("I3").Value = 1042.5714
$ws.Range("J3").Value = 11832.25
$ws.Range("K3").Value = 3127.7142
$ws.Range("L3").Value = 35496.75
$ws.Range("M3").Value = -3015.7142
$ws.Range("N3").Value = -35720.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 350963.44
$ws.Range("I122").Value = 791312.1
$ws.Range("K122").Value = 2373936.3
$ws.Range("M122").Value = -2371486.3
$ws.Range("H123").Value = 53618
$ws.Range("J123").Value = 53618
$ws.Range("L123").Value = 53618
$ws.Range("N123").Value = -58518
$ws.Range("H132").Value = 4043.647
$ws.Range("I132").Value = 3601.25
$ws.Range("J132").Value = 6108.1665
$ws.Range("K132").Value = 10803.75
$ws.Range("L132").Value = 18324.4995
$ws.Range("M132").Value = -8273.75
$ws.Range("N132").Value = -23384.4995
$ws.Range("H139").Value = 90000
$ws.Range("J139").Value = 90000
$ws.Range("L139").Value = 90000
$ws.Range("N139").Value = -100280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2356.3928
$ws.Range("I16").Value = 1695.5
$ws.Range("J16").Value = 4779.6665
$ws.Range("K16").Value = 1695.5
$ws.Range("L16").Value = 4779.6665
$ws.Range("M16").Value = -1525.5
$ws.Range("N16").Value = -5119.6665
$ws.Range("H122").Value = 14194.177
$ws.Range("I122").Value = 16672
$ws.Range("K122").Value = 50016
$ws.Range("M122").Value = -47566
$ws.Range("H132").Value = 4151.8687
$ws.Range("I132").Value = 3918.1904
$ws.Range("K132").Value = 11754.5712
$ws.Range("M132").Value = -9224.5712

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1502567.6
$ws.Range("I81").Value = 2610126.8
$ws.Range("J81").Value = 25822.166
$ws.Range("K81").Value = 5220253.6
$ws.Range("L81").Value = 51644.332
$ws.Range("M81").Value = -5219192.6
$ws.Range("N81").Value = -53766.332
$ws.Range("H84").Value = 1502567.6
$ws.Range("I84").Value = 2610126.8
$ws.Range("J84").Value = 25822.166
$ws.Range("K84").Value = 26101268
$ws.Range("L84").Value = 258221.66
$ws.Range("M84").Value = -26095964
$ws.Range("N84").Value = -268829.66
$ws.Range("H100").Value = 557943.25
$ws.Range("I100").Value = 751347.2
$ws.Range("K100").Value = 1502694.4
$ws.Range("M100").Value = -1502153.4
$ws.Range("H111").Value = 92500
$ws.Range("J111").Value = 92500
$ws.Range("L111").Value = 92500
$ws.Range("N111").Value = -100680
$ws.Range("I132").Value = 874.7273
$ws.Range("K132").Value = 2624.1819
$ws.Range("M132").Value = -94.18190000000004
$ws.Range("H136").Value = 7234.1
$ws.Range("I136").Value = 3283.348
$ws.Range("J136").Value = 10599.556
$ws.Range("K136").Value = 9850.044
$ws.Range("L136").Value = 31798.668
$ws.Range("M136").Value = -7300.044
$ws.Range("N136").Value = -36898.66800000001
